$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 113: Tue 16 Jan 2024 - Rest, with a logged time/intensity entry
$ws.Range("A108:G108").Copy($ws.Range("A113:G113"))
$ws.Range("A113").Value = 45307
$ws.Range("B113").Value = "Rest"
$ws.Range("C113").Value = "Injured"
$ws.Range("D113").Value = 0
$ws.Range("E113").Clear()
$ws.Range("F113").Value = 0.625
$ws.Range("G113").Value = 7

# Row 114: Wed 17 Jan 2024 - Cardio
$ws.Range("A112:D112").Copy($ws.Range("A114:D114"))
$ws.Range("A114").Value = 45308
$ws.Range("B114").Value = "Cardio"
$ws.Range("C114").Value = "Injured"
$ws.Range("D114").Value = 15

# Row 115: Thu 18 Jan 2024 - Rest
$ws.Range("A112:D112").Copy($ws.Range("A115:D115"))
$ws.Range("A115").Value = 45309
$ws.Range("B115").Value = "Rest"
$ws.Range("C115").Value = "Injured"
$ws.Range("D115").Value = 0

# Row 116: Fri 19 Jan 2024 - Rest
$ws.Range("A112:D112").Copy($ws.Range("A116:D116"))
$ws.Range("A116").Value = 45310
$ws.Range("B116").Value = "Rest"
$ws.Range("C116").Value = "Injured"
$ws.Range("D116").Value = 0

# Row 117: Sat 20 Jan 2024 - Cardio, with logged time/intensity
$ws.Range("A108:G108").Copy($ws.Range("A117:G117"))
$ws.Range("A117").Value = 45311
$ws.Range("B117").Value = "Cardio"
$ws.Range("C117").Value = "Injured"
$ws.Range("D117").Value = 15
$ws.Range("E117").Clear()
$ws.Range("F117").Value = 0.60416666666666663
$ws.Range("G117").Value = 7

# Row 118: Sun 21 Jan 2024 - Rest
$ws.Range("A112:D112").Copy($ws.Range("A118:D118"))
$ws.Range("A118").Value = 45312
$ws.Range("B118").Value = "Rest"
$ws.Range("C118").Value = "Injured"
$ws.Range("D118").Value = 0

# Row 119: Mon 22 Jan 2024 - Rest
$ws.Range("A112:D112").Copy($ws.Range("A119:D119"))
$ws.Range("A119").Value = 45313
$ws.Range("B119").Value = "Rest"
$ws.Range("C119").Value = "Injured"
$ws.Range("D119").Value = 0

# Row 120: Tue 23 Jan 2024 - Cardio, with logged time only (no intensity)
$ws.Range("A108:F108").Copy($ws.Range("A120:F120"))
$ws.Range("A120").Value = 45314
$ws.Range("B120").Value = "Cardio"
$ws.Range("C120").Value = "Injured"
$ws.Range("D120").Value = 15
$ws.Range("E120").Clear()
$ws.Range("F120").Value = 0.58333333333333337

# Row 121: Wed 24 Jan 2024 - Rest
$ws.Range("A112:D112").Copy($ws.Range("A121:D121"))
$ws.Range("A121").Value = 45315
$ws.Range("B121").Value = "Rest"
$ws.Range("C121").Value = "Injured"
$ws.Range("D121").Value = 0

# Row 122: Thu 25 Jan 2024 - Rest
$ws.Range("A112:D112").Copy($ws.Range("A122:D122"))
$ws.Range("A122").Value = 45316
$ws.Range("B122").Value = "Rest"
$ws.Range("C122").Value = "Injured"
$ws.Range("D122").Value = 0

# Row 123: Fri 26 Jan 2024 - Cardio, with logged time/intensity
$ws.Range("A108:G108").Copy($ws.Range("A123:G123"))
$ws.Range("A123").Value = 45317
$ws.Range("B123").Value = "Cardio"
$ws.Range("C123").Value = "Injured"
$ws.Range("D123").Value = 15
$ws.Range("E123").Clear()
$ws.Range("F123").Value = 0.54166666666666663
$ws.Range("G123").Value = 7

# Row 124: Sat 27 Jan 2024 - Rest
$ws.Range("A112:D112").Copy($ws.Range("A124:D124"))
$ws.Range("A124").Value = 45318
$ws.Range("B124").Value = "Rest"
$ws.Range("C124").Value = "Injured"
$ws.Range("D124").Value = 0

# Row 125: Sun 28 Jan 2024 - Rest
$ws.Range("A112:D112").Copy($ws.Range("A125:D125"))
$ws.Range("A125").Value = 45319
$ws.Range("B125").Value = "Rest"
$ws.Range("C125").Value = "Injured"
$ws.Range("D125").Value = 0

# Row 126: Mon 29 Jan 2024 - Rest
$ws.Range("A112:D112").Copy($ws.Range("A126:D126"))
$ws.Range("A126").Value = 45320
$ws.Range("B126").Value = "Rest"
$ws.Range("C126").Value = "Injured"
$ws.Range("D126").Value = 0

# Row 127: Tue 30 Jan 2024 - Rest
$ws.Range("A112:D112").Copy($ws.Range("A127:D127"))
$ws.Range("A127").Value = 45321
$ws.Range("B127").Value = "Rest"
$ws.Range("C127").Value = "Injured"
$ws.Range("D127").Value = 0

# Restore the selection to reflect the last-edited cell among the new rows
[void]$ws.Range("G123").Select()
